$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value2 = "ECs"
$ws.Cells.Item(2, 2).Value2 = "Inhba"
$ws.Cells.Item(2, 3).Value2 = "Acvr1"
$ws.Cells.Item(2, 4).Value2 = "ECs"
$ws.Cells.Item(2, 5).Value2 = 2
$ws.Cells.Item(2, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(2, 7).Value2 = 3.675031333333333
$ws.Cells.Item(2, 8).Value2 = 11.025094
$ws.Cells.Item(2, 9).Value2 = 0.2032371147293133
$ws.Cells.Item(2, 10).Value2 = 0.2032371147293133
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 12).Value2 = 1
$ws.Cells.Item(2, 13).Value2 = 4.695610666666666
$ws.Cells.Item(2, 14).Value2 = 14.086832
$ws.Cells.Item(2, 15).Value2 = 0.1802066564018305
$ws.Cells.Item(2, 16).Value2 = 0.1802066564018305
$ws.Cells.Item(2, 17).Value2 = 17.25651632913422
$ws.Cells.Item(2, 18).Value2 = 155.308646962208
$ws.Cells.Item(2, 19).Value2 = 0.03662468090212477
$ws.Cells.Item(2, 20).Value2 = 0.03662468090212476

# Row 3
$ws.Cells.Item(3, 1).Value2 = "ECs"
$ws.Cells.Item(3, 2).Value2 = "Inhba"
$ws.Cells.Item(3, 3).Value2 = "Acvr1"
$ws.Cells.Item(3, 4).Value2 = "FAPs"
$ws.Cells.Item(3, 5).Value2 = 2
$ws.Cells.Item(3, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(3, 7).Value2 = 3.675031333333333
$ws.Cells.Item(3, 8).Value2 = 11.025094
$ws.Cells.Item(3, 9).Value2 = 0.2032371147293133
$ws.Cells.Item(3, 10).Value2 = 0.2032371147293133
$ws.Cells.Item(3, 11).Value2 = 3
$ws.Cells.Item(3, 12).Value2 = 1
$ws.Cells.Item(3, 13).Value2 = 15.51448033333333
$ws.Cells.Item(3, 14).Value2 = 46.543441
$ws.Cells.Item(3, 15).Value2 = 0.5954098039960916
$ws.Cells.Item(3, 16).Value2 = 0.5954098039960916
$ws.Cells.Item(3, 17).Value2 = 57.01620134538377
$ws.Cells.Item(3, 18).Value2 = 513.145812108454
$ws.Cells.Item(3, 19).Value2 = 0.1210093706457116
$ws.Cells.Item(3, 20).Value2 = 0.1210093706457116

# Row 4
$ws.Cells.Item(4, 1).Value2 = "ECs"
$ws.Cells.Item(4, 2).Value2 = "Inhba"
$ws.Cells.Item(4, 3).Value2 = "Acvr1"
$ws.Cells.Item(4, 4).Value2 = "sCs"
$ws.Cells.Item(4, 5).Value2 = 2
$ws.Cells.Item(4, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(4, 7).Value2 = 3.675031333333333
$ws.Cells.Item(4, 8).Value2 = 11.025094
$ws.Cells.Item(4, 9).Value2 = 0.2032371147293133
$ws.Cells.Item(4, 10).Value2 = 0.2032371147293133
$ws.Cells.Item(4, 11).Value2 = 3
$ws.Cells.Item(4, 12).Value2 = 1
$ws.Cells.Item(4, 13).Value2 = 5.846719333333333
$ws.Cells.Item(4, 14).Value2 = 17.540158
$ws.Cells.Item(4, 15).Value2 = 0.2243835396020779
$ws.Cells.Item(4, 16).Value2 = 0.2243835396020779
$ws.Cells.Item(4, 17).Value2 = 21.48687674720577
$ws.Cells.Item(4, 18).Value2 = 193.381890724852
$ws.Cells.Item(4, 19).Value2 = 0.04560306318147692
$ws.Cells.Item(4, 20).Value2 = 0.04560306318147692

# Row 5
$ws.Cells.Item(5, 1).Value2 = "FAPs"
$ws.Cells.Item(5, 2).Value2 = "Inhba"
$ws.Cells.Item(5, 3).Value2 = "Acvr1"
$ws.Cells.Item(5, 4).Value2 = "ECs"
$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 6).Value2 = 1
$ws.Cells.Item(5, 7).Value2 = 10.108494
$ws.Cells.Item(5, 8).Value2 = 30.325482
$ws.Cells.Item(5, 9).Value2 = 0.5590213983169419
$ws.Cells.Item(5, 10).Value2 = 0.5590213983169419
$ws.Cells.Item(5, 11).Value2 = 3
$ws.Cells.Item(5, 12).Value2 = 1
$ws.Cells.Item(5, 13).Value2 = 4.695610666666666
$ws.Cells.Item(5, 14).Value2 = 14.086832
$ws.Cells.Item(5, 15).Value2 = 0.1802066564018305
$ws.Cells.Item(5, 16).Value2 = 0.1802066564018305
$ws.Cells.Item(5, 17).Value2 = 47.465552250336
$ws.Cells.Item(5, 18).Value2 = 427.189970253024
$ws.Cells.Item(5, 19).Value2 = 0.100739377047772
$ws.Cells.Item(5, 20).Value2 = 0.100739377047772

# Row 6
$ws.Cells.Item(6, 1).Value2 = "FAPs"
$ws.Cells.Item(6, 2).Value2 = "Inhba"
$ws.Cells.Item(6, 3).Value2 = "Acvr1"
$ws.Cells.Item(6, 4).Value2 = "FAPs"
$ws.Cells.Item(6, 5).Value2 = 3
$ws.Cells.Item(6, 6).Value2 = 1
$ws.Cells.Item(6, 7).Value2 = 10.108494
$ws.Cells.Item(6, 8).Value2 = 30.325482
$ws.Cells.Item(6, 9).Value2 = 0.5590213983169419
$ws.Cells.Item(6, 10).Value2 = 0.5590213983169419
$ws.Cells.Item(6, 11).Value2 = 3
$ws.Cells.Item(6, 12).Value2 = 1
$ws.Cells.Item(6, 13).Value2 = 15.51448033333333
$ws.Cells.Item(6, 14).Value2 = 46.543441
$ws.Cells.Item(6, 15).Value2 = 0.5954098039960916
$ws.Cells.Item(6, 16).Value2 = 0.5954098039960916
$ws.Cells.Item(6, 17).Value2 = 156.828031362618
$ws.Cells.Item(6, 18).Value2 = 1411.452282263562
$ws.Cells.Item(6, 19).Value2 = 0.3328468212015114
$ws.Cells.Item(6, 20).Value2 = 0.3328468212015114

# Row 7
$ws.Cells.Item(7, 1).Value2 = "FAPs"
$ws.Cells.Item(7, 2).Value2 = "Inhba"
$ws.Cells.Item(7, 3).Value2 = "Acvr1"
$ws.Cells.Item(7, 4).Value2 = "sCs"
$ws.Cells.Item(7, 5).Value2 = 3
$ws.Cells.Item(7, 6).Value2 = 1
$ws.Cells.Item(7, 7).Value2 = 10.108494
$ws.Cells.Item(7, 8).Value2 = 30.325482
$ws.Cells.Item(7, 9).Value2 = 0.5590213983169419
$ws.Cells.Item(7, 10).Value2 = 0.5590213983169419
$ws.Cells.Item(7, 11).Value2 = 3
$ws.Cells.Item(7, 12).Value2 = 1
$ws.Cells.Item(7, 13).Value2 = 5.846719333333333
$ws.Cells.Item(7, 14).Value2 = 17.540158
$ws.Cells.Item(7, 15).Value2 = 0.2243835396020779
$ws.Cells.Item(7, 16).Value2 = 0.2243835396020779
$ws.Cells.Item(7, 17).Value2 = 59.101527300684
$ws.Cells.Item(7, 18).Value2 = 531.913745706156
$ws.Cells.Item(7, 19).Value2 = 0.1254352000676585
$ws.Cells.Item(7, 20).Value2 = 0.1254352000676585

# Row 8
$ws.Cells.Item(8, 1).Value2 = "sCs"
$ws.Cells.Item(8, 2).Value2 = "Inhba"
$ws.Cells.Item(8, 3).Value2 = "Acvr1"
$ws.Cells.Item(8, 4).Value2 = "ECs"
$ws.Cells.Item(8, 5).Value2 = 3
$ws.Cells.Item(8, 6).Value2 = 1
$ws.Cells.Item(8, 7).Value2 = 4.298956
$ws.Cells.Item(8, 8).Value2 = 12.896868
$ws.Cells.Item(8, 9).Value2 = 0.2377414869537448
$ws.Cells.Item(8, 10).Value2 = 0.2377414869537448
$ws.Cells.Item(8, 11).Value2 = 3
$ws.Cells.Item(8, 12).Value2 = 1
$ws.Cells.Item(8, 13).Value2 = 4.695610666666666
$ws.Cells.Item(8, 14).Value2 = 14.086832
$ws.Cells.Item(8, 15).Value2 = 0.1802066564018305
$ws.Cells.Item(8, 16).Value2 = 0.1802066564018305
$ws.Cells.Item(8, 17).Value2 = 20.18622364913067
$ws.Cells.Item(8, 18).Value2 = 181.676012842176
$ws.Cells.Item(8, 19).Value2 = 0.04284259845193376
$ws.Cells.Item(8, 20).Value2 = 0.04284259845193375

# Row 9
$ws.Cells.Item(9, 1).Value2 = "sCs"
$ws.Cells.Item(9, 2).Value2 = "Inhba"
$ws.Cells.Item(9, 3).Value2 = "Acvr1"
$ws.Cells.Item(9, 4).Value2 = "FAPs"
$ws.Cells.Item(9, 5).Value2 = 3
$ws.Cells.Item(9, 6).Value2 = 1
$ws.Cells.Item(9, 7).Value2 = 4.298956
$ws.Cells.Item(9, 8).Value2 = 12.896868
$ws.Cells.Item(9, 9).Value2 = 0.2377414869537448
$ws.Cells.Item(9, 10).Value2 = 0.2377414869537448
$ws.Cells.Item(9, 11).Value2 = 3
$ws.Cells.Item(9, 12).Value2 = 1
$ws.Cells.Item(9, 13).Value2 = 15.51448033333333
$ws.Cells.Item(9, 14).Value2 = 46.543441
$ws.Cells.Item(9, 15).Value2 = 0.5954098039960916
$ws.Cells.Item(9, 16).Value2 = 0.5954098039960916
$ws.Cells.Item(9, 17).Value2 = 66.69606831586535
$ws.Cells.Item(9, 18).Value2 = 600.2646148427881
$ws.Cells.Item(9, 19).Value2 = 0.1415536121488685
$ws.Cells.Item(9, 20).Value2 = 0.1415536121488685

# Row 10
$ws.Cells.Item(10, 1).Value2 = "sCs"
$ws.Cells.Item(10, 2).Value2 = "Inhba"
$ws.Cells.Item(10, 3).Value2 = "Acvr1"
$ws.Cells.Item(10, 4).Value2 = "sCs"
$ws.Cells.Item(10, 5).Value2 = 3
$ws.Cells.Item(10, 6).Value2 = 1
$ws.Cells.Item(10, 7).Value2 = 4.298956
$ws.Cells.Item(10, 8).Value2 = 12.896868
$ws.Cells.Item(10, 9).Value2 = 0.2377414869537448
$ws.Cells.Item(10, 10).Value2 = 0.2377414869537448
$ws.Cells.Item(10, 11).Value2 = 3
$ws.Cells.Item(10, 12).Value2 = 1
$ws.Cells.Item(10, 13).Value2 = 5.846719333333333
$ws.Cells.Item(10, 14).Value2 = 17.540158
$ws.Cells.Item(10, 15).Value2 = 0.2243835396020779
$ws.Cells.Item(10, 16).Value2 = 0.2243835396020779
$ws.Cells.Item(10, 17).Value2 = 25.13478915834934
$ws.Cells.Item(10, 18).Value2 = 226.213102425144
$ws.Cells.Item(10, 19).Value2 = 0.05334527635294248
$ws.Cells.Item(10, 20).Value2 = 0.05334527635294247
